$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are plain decimal numbers (e.g. "578.19"). Force those
# specific cells to Text format first so Excel keeps them as text (matching the rest
# of the column, which is all inline/shared text) instead of silently converting them
# to floating point numbers.
$textCells = @("D5", "D6", "D8", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D30", "D32", "D34", "D35", "D36", "D38", "D40", "D42", "D44", "D46", "D47", "D48")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.602.58"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "3.255.29"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "578.19"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "181.59"
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "3.254.62"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("E12").Value = "  +4.98%  "
$ws.Range("D13").Value = "3.818.53"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("D14").Value = "0.137"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "28.79"
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").Value = "67.574.70"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "3.257.45"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").Value = "5.85"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "13.56"
$ws.Range("E20").Value = "  +5.41%  "
$ws.Range("D21").Value = "377.07"
$ws.Range("E21").Value = "  +5.82%  "
$ws.Range("D22").Value = "7.64"
$ws.Range("E22").Value = "  +5.19%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "71.46"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "5.78"
$ws.Range("E30").Value = "  +6.87%  "
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").Value = "22.69"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  +5.56%  "
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").Value = "163.85"
$ws.Range("E36").Value = "  +6.38%  "
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("D38").Value = "0.851"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("E39").Value = "  +5.18%  "
$ws.Range("D40").Value = "26.88"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("E41").Value = "  +12.37%  "
$ws.Range("D42").Value = "4.58"
$ws.Range("E42").Value = "  +9.23%  "
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("D44").Value = "364.98"
$ws.Range("E44").Value = "  +10.74%  "
$ws.Range("D45").Value = "2.737.96"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").Value = "25.44"
$ws.Range("E46").Value = "  +4.13%  "
$ws.Range("D47").Value = "40.90"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").Value = "0.0684"
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  +6.72%  "
$ws.Range("E51").Value = "  +0.32%  "
